# Update the "as of" date in the confidential disclosure footer and refresh
# the Weight / Percent Change figures in the holdings table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; temporarily unprotect so values can be edited,
# then restore protection at the end.
$ws.Unprotect()

# --- Update the "as of" date in the confidential disclosure text (cell A18) ---
# Use Range.Replace (text find/replace) rather than re-assigning the whole
# multi-line string, then AutoFit the row back so no spurious custom row
# height gets baked in.
$ws.Range("A18").Replace("2021-03-30", "2021-03-31") | Out-Null
$ws.Rows.Item(18).AutoFit()

# --- Update Weight (D) and Percent Change (E) values for rows 2-15 ---
$values = @(
    @{ Row = 2;  D = 0.05587839703200414;  E = 0.006084501557632294 }
    @{ Row = 3;  D = 0.02344246176107796;  E = -0.002842639593908736 }
    @{ Row = 4;  D = 0.03171311186695153;  E = 0.004438440756464601 }
    @{ Row = 5;  D = 0.03204086434175471;  E = -0.006048387096774244 }
    @{ Row = 6;  D = 0.03404273837872682;  E = 0.009827420901246553 }
    @{ Row = 7;  D = 0.01906812219591438;  E = -0.001497753369945176 }
    @{ Row = 8;  D = 0.004791849979124478; E = 0.03448275862068972 }
    @{ Row = 9;  D = 0.006677446685824093; E = 0.006109979633401208 }
    @{ Row = 10; D = 0.06821603375986486;  E = 0.01555023923444998 }
    @{ Row = 11; D = 0.06829763188637188;  E = 0.0155316606929512 }
    @{ Row = 12; D = 0.1481930774869409;   E = -0.005579619704867644 }
    @{ Row = 13; D = 0.3921871153838308;   E = -0.001414802369793966 }
    @{ Row = 14; D = 0.1154511492416134;   E = -0.001236858379715433 }
    @{ Row = 15; D = 0.9999999999999999;   E = 0.001329369477676812 }
)

foreach ($entry in $values) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value2 = $entry.D
    $ws.Cells.Item($r, 5).Value2 = $entry.E
}

# Restore sheet protection to its original state.
$ws.Protect()
